# Estandarización de las columnas "nombre" por "nom" en la tabla usuario.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("usuario")

# Renombrar los encabezados de columna F2 y G2.
$ws.Range("F2").Value = "nom"
$ws.Range("G2").Value = "nom_preferido"

# Forzar el recalculo para que las formulas de la columna J (que generan
# los "insert into ...") reflejen el nuevo nombre de columna.
$excel.CalculateFullRebuild()

# Dejar seleccionada la celda F23, tal como quedo en el archivo final.
$ws.Range("F23").Select()

$wb.Save()
